$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the bordered style from B3 onto the new B13:B24 cells
$ws.Cells.Item(3, 2).Copy()
$ws.Range("B13:B24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Continue the numeric series 11..20 down column B (rows 13-22)
for ($i = 0; $i -lt 10; $i++) {
    $row = 13 + $i
    $value = 11 + $i
    $ws.Cells.Item($row, 2).Value = $value
}

# Move the two average formulas from row 13 down to row 23
$ws.Range("E23").Formula = "=SUM(E3:E12)/COUNT(E3:E12)"
$ws.Range("F23").Formula = "=SUM(F3:F12)/COUNT(F3:F12)"
$ws.Range("E13:F13").ClearContents()

# Update the active selection to match
$ws.Range("E23").Select()
